$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vegfc"
$ws.Range("C2").Value = "Kdr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 3.911257666666666
$ws.Range("H2").Value = 11.733773
$ws.Range("I2").Value = 0.4115343446855154
$ws.Range("J2").Value = 0.4115343446855154
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 162.98837
$ws.Range("N2").Value = 488.96511
$ws.Range("O2").Value = 0.9909539753179891
$ws.Range("P2").Value = 0.9909539753179891
$ws.Range("Q2").Value = 637.4895117400033
$ws.Range("R2").Value = 5737.405605660029
$ws.Range("S2").Value = 0.407811594845995
$ws.Range("T2").Value = 0.407811594845995

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vegfc"
$ws.Range("C3").Value = "Kdr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 3.911257666666666
$ws.Range("H3").Value = 11.733773
$ws.Range("I3").Value = 0.4115343446855154
$ws.Range("J3").Value = 0.4115343446855154
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.6513563333333333
$ws.Range("N3").Value = 1.954069
$ws.Range("O3").Value = 0.003960185305646138
$ws.Range("P3").Value = 0.003960185305646138
$ws.Range("Q3").Value = 2.547622452481888
$ws.Range("R3").Value = 22.928602072337
$ws.Range("S3").Value = 0.001629752264592291
$ws.Range("T3").Value = 0.001629752264592291

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Vegfc"
$ws.Range("C4").Value = "Kdr"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 3.911257666666666
$ws.Range("H4").Value = 11.733773
$ws.Range("I4").Value = 0.4115343446855154
$ws.Range("J4").Value = 0.4115343446855154
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.8364996666666666
$ws.Range("N4").Value = 2.509499
$ws.Range("O4").Value = 0.005085839376364744
$ws.Range("P4").Value = 0.005085839376364744
$ws.Range("Q4").Value = 3.271765734414111
$ws.Range("R4").Value = 29.445891609727
$ws.Range("S4").Value = 0.002092997574928055
$ws.Range("T4").Value = 0.002092997574928055

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Vegfc"
$ws.Range("C5").Value = "Kdr"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.100181333333333
$ws.Range("H5").Value = 12.300544
$ws.Range("I5").Value = 0.4314124974392592
$ws.Range("J5").Value = 0.4314124974392592
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 162.98837
$ws.Range("N5").Value = 488.96511
$ws.Range("O5").Value = 0.9909539753179891
$ws.Range("P5").Value = 0.9909539753179891
$ws.Range("Q5").Value = 668.2818722244267
$ws.Range("R5").Value = 6014.53685001984
$ws.Range("S5").Value = 0.4275099293392957
$ws.Range("T5").Value = 0.4275099293392957

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Vegfc"
$ws.Range("C6").Value = "Kdr"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.100181333333333
$ws.Range("H6").Value = 12.300544
$ws.Range("I6").Value = 0.4314124974392592
$ws.Range("J6").Value = 0.4314124974392592
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.6513563333333333
$ws.Range("N6").Value = 1.954069
$ws.Range("O6").Value = 0.003960185305646138
$ws.Range("P6").Value = 0.003960185305646138
$ws.Range("Q6").Value = 2.670679079281778
$ws.Range("R6").Value = 24.036111713536
$ws.Range("S6").Value = 0.001708473433031057
$ws.Range("T6").Value = 0.001708473433031057

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Vegfc"
$ws.Range("C7").Value = "Kdr"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.100181333333333
$ws.Range("H7").Value = 12.300544
$ws.Range("I7").Value = 0.4314124974392592
$ws.Range("J7").Value = 0.4314124974392592
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.8364996666666666
$ws.Range("N7").Value = 2.509499
$ws.Range("O7").Value = 0.005085839376364744
$ws.Range("P7").Value = 0.005085839376364744
$ws.Range("Q7").Value = 3.429800318606222
$ws.Range("R7").Value = 30.868202867456
$ws.Range("S7").Value = 0.002194094666932439
$ws.Range("T7").Value = 0.002194094666932439

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Vegfc"
$ws.Range("C8").Value = "Kdr"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.492646666666667
$ws.Range("H8").Value = 4.47794
$ws.Range("I8").Value = 0.1570531578752254
$ws.Range("J8").Value = 0.1570531578752254
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 162.98837
$ws.Range("N8").Value = 488.96511
$ws.Range("O8").Value = 0.9909539753179891
$ws.Range("P8").Value = 0.9909539753179891
$ws.Range("Q8").Value = 243.2840471859333
$ws.Range("R8").Value = 2189.5564246734
$ws.Range("S8").Value = 0.1556324511326984
$ws.Range("T8").Value = 0.1556324511326984

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Vegfc"
$ws.Range("C9").Value = "Kdr"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.492646666666667
$ws.Range("H9").Value = 4.47794
$ws.Range("I9").Value = 0.1570531578752254
$ws.Range("J9").Value = 0.1570531578752254
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6513563333333333
$ws.Range("N9").Value = 1.954069
$ws.Range("O9").Value = 0.003960185305646138
$ws.Range("P9").Value = 0.003960185305646138
$ws.Range("Q9").Value = 0.9722448597622222
$ws.Range("R9").Value = 8.750203737860001
$ws.Range("S9").Value = 0.0006219596080227907
$ws.Range("T9").Value = 0.0006219596080227907

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Vegfc"
$ws.Range("C10").Value = "Kdr"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.492646666666667
$ws.Range("H10").Value = 4.47794
$ws.Range("I10").Value = 0.1570531578752254
$ws.Range("J10").Value = 0.1570531578752254
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.8364996666666666
$ws.Range("N10").Value = 2.509499
$ws.Range("O10").Value = 0.005085839376364744
$ws.Range("P10").Value = 0.005085839376364744
$ws.Range("Q10").Value = 1.248598439117778
$ws.Range("R10").Value = 11.23738595206
$ws.Range("S10").Value = 0.00079874713450425
$ws.Range("T10").Value = 0.00079874713450425
